# Applies the scheduled-runner profit recalculation update to the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (Tonberry_Profits workbook).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (92 cell(s)) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2670.625
$ws.Range("I6").Value = 1149
$ws.Range("K6").Value = 3447
$ws.Range("M6").Value = -3335
$ws.Range("H18").Value = 15359.2
$ws.Range("I18").Value = 898
$ws.Range("K18").Value = 898
$ws.Range("M18").Value = -614
$ws.Range("H43").Value = 1613.7273
$ws.Range("I43").Value = 1625.1
$ws.Range("K43").Value = 1625.1
$ws.Range("M43").Value = -1556.1
$ws.Range("H51").Value = 4899.7
$ws.Range("J51").Value = 5110.8887
$ws.Range("L51").Value = 5110.8887
$ws.Range("N51").Value = -6078.8887
$ws.Range("H70").Value = 1045
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 1154
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 3462
$ws.Range("M70").Value = -1230
$ws.Range("N70").Value = -4002
$ws.Range("H73").Value = 1045
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 1154
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 3462
$ws.Range("M73").Value = -564
$ws.Range("N73").Value = -5334
$ws.Range("H76").Value = 4651.5713
$ws.Range("I76").Value = 4427
$ws.Range("J76").Value = 5999
$ws.Range("K76").Value = 4427
$ws.Range("L76").Value = 5999
$ws.Range("M76").Value = -4112
$ws.Range("N76").Value = -6629
$ws.Range("H79").Value = 4651.5713
$ws.Range("I79").Value = 4427
$ws.Range("J79").Value = 5999
$ws.Range("K79").Value = 4427
$ws.Range("L79").Value = 5999
$ws.Range("M79").Value = -3335
$ws.Range("N79").Value = -8183
$ws.Range("H92").Value = 17857388
$ws.Range("I92").Value = 20833530
$ws.Range("J92").Value = 535.5
$ws.Range("K92").Value = 20833530
$ws.Range("L92").Value = 535.5
$ws.Range("M92").Value = -20832282
$ws.Range("N92").Value = -3031.5
$ws.Range("H106").Value = 2990.1924
$ws.Range("I106").Value = 2281.45
$ws.Range("K106").Value = 2281.45
$ws.Range("M106").Value = -1650.45
$ws.Range("H107").Value = 449.15384
$ws.Range("I107").Value = 511.63635
$ws.Range("K107").Value = 511.63635
$ws.Range("M107").Value = 1408.36365
$ws.Range("H112").Value = 2191.0222
$ws.Range("J112").Value = 2191.0222
$ws.Range("L112").Value = 6573.0666
$ws.Range("N112").Value = -8789.0666
$ws.Range("H125").Value = 833
$ws.Range("I125").Value = 1099.5
$ws.Range("J125").Value = 300
$ws.Range("K125").Value = 9895.5
$ws.Range("L125").Value = 2700
$ws.Range("M125").Value = -7435.5
$ws.Range("N125").Value = -7620
$ws.Range("H132").Value = 1222.6072
$ws.Range("I132").Value = 1177.7916
$ws.Range("K132").Value = 3533.3748
$ws.Range("M132").Value = -1003.3748
$ws.Range("H137").Value = 49342.24
$ws.Range("I137").Value = 1257.8
$ws.Range("J137").Value = 93055.37
$ws.Range("K137").Value = 3773.4
$ws.Range("L137").Value = 279166.11
$ws.Range("M137").Value = -1223.4
$ws.Range("N137").Value = -284266.11
$ws.Range("H138").Value = 3633.16
$ws.Range("J138").Value = 2949.3
$ws.Range("L138").Value = 8847.900000000001
$ws.Range("N138").Value = -19127.9
$ws.Range("H141").Value = 1402113.9
$ws.Range("I141").Value = 1867885.2
$ws.Range("J141").Value = 4799.8
$ws.Range("K141").Value = 5603655.6
$ws.Range("L141").Value = 14399.4
$ws.Range("M141").Value = -5598475.6
$ws.Range("N141").Value = -24759.4

# ---- Sheet: ARM (27 cell(s)) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10673.695
$ws.Range("I32").Value = 8312.775
$ws.Range("J32").Value = 26413.166
$ws.Range("K32").Value = 8312.775
$ws.Range("L32").Value = 26413.166
$ws.Range("M32").Value = -8025.775
$ws.Range("N32").Value = -26987.166
$ws.Range("H45").Value = 5001273
$ws.Range("I45").Value = 10000722
$ws.Range("K45").Value = 10000722
$ws.Range("M45").Value = -10000345
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50636
$ws.Range("H76").Value = 25600
$ws.Range("J76").Value = 25600
$ws.Range("L76").Value = 25600
$ws.Range("N76").Value = -26276
$ws.Range("H79").Value = 25600
$ws.Range("J79").Value = 25600
$ws.Range("L79").Value = 25600
$ws.Range("N79").Value = -27940
$ws.Range("H122").Value = 22577.055
$ws.Range("I122").Value = 25149.25
$ws.Range("K122").Value = 75447.75
$ws.Range("M122").Value = -72997.75

# ---- Sheet: BSM (25 cell(s)) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2775
$ws.Range("I20").Value = 2033.3334
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 2033.3334
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -1786.3334
$ws.Range("N20").Value = -5494
$ws.Range("H99").Value = 1011
$ws.Range("I99").Value = 1011
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1011
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 487
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 3239.889
$ws.Range("I107").Value = 3239.889
$ws.Range("K107").Value = 3239.889
$ws.Range("M107").Value = -1319.889
$ws.Range("H134").Value = 4906.028
$ws.Range("I134").Value = 5332.7407
$ws.Range("J134").Value = 3625.889
$ws.Range("K134").Value = 15998.2221
$ws.Range("L134").Value = 10877.667
$ws.Range("M134").Value = -13463.2221
$ws.Range("N134").Value = -15947.667

# ---- Sheet: CRP (25 cell(s)) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2526.2812
$ws.Range("I31").Value = 2163.158
$ws.Range("J31").Value = 3057
$ws.Range("K31").Value = 2163.158
$ws.Range("L31").Value = 3057
$ws.Range("M31").Value = -1868.158
$ws.Range("N31").Value = -3647
$ws.Range("H34").Value = 2526.2812
$ws.Range("I34").Value = 2163.158
$ws.Range("J34").Value = 3057
$ws.Range("K34").Value = 2163.158
$ws.Range("L34").Value = 3057
$ws.Range("M34").Value = -1961.158
$ws.Range("N34").Value = -3461
$ws.Range("H98").Value = 69993
$ws.Range("J98").Value = 69993
$ws.Range("L98").Value = 69993
$ws.Range("N98").Value = -74485
$ws.Range("H132").Value = 1465.4839
$ws.Range("I132").Value = 892.5
$ws.Range("J132").Value = 3430
$ws.Range("K132").Value = 2677.5
$ws.Range("L132").Value = 10290
$ws.Range("M132").Value = -147.5
$ws.Range("N132").Value = -15350

# ---- Sheet: CUL (23 cell(s)) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6504.077
$ws.Range("I56").Value = 6504.077
$ws.Range("K56").Value = 6504.077
$ws.Range("M56").Value = -5974.077
$ws.Range("H68").Value = 885.625
$ws.Range("I68").Value = 797.8570999999999
$ws.Range("K68").Value = 2393.5713
$ws.Range("M68").Value = -1582.5713
$ws.Range("H71").Value = 885.625
$ws.Range("I71").Value = 797.8570999999999
$ws.Range("K71").Value = 7180.7139
$ws.Range("M71").Value = -3124.7139
$ws.Range("H110").Value = 4390.5835
$ws.Range("I110").Value = 2562.3333
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 7686.999899999999
$ws.Range("L110").Value = 15000
$ws.Range("M110").Value = -3596.999899999999
$ws.Range("N110").Value = -23180
$ws.Range("H131").Value = 25989.875
$ws.Range("J131").Value = 33142.4
$ws.Range("L131").Value = 99427.20000000001
$ws.Range("N131").Value = -109507.2

# ---- Sheet: GSM (8 cell(s)) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 26700
$ws.Range("J93").Value = 26700
$ws.Range("L93").Value = 26700
$ws.Range("N93").Value = -30444
$ws.Range("H97").Value = 3256
$ws.Range("I97").Value = 3256
$ws.Range("K97").Value = 3256
$ws.Range("M97").Value = -2760

# ---- Sheet: LTW (26 cell(s)) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1962.0714
$ws.Range("I46").Value = 1611.4286
$ws.Range("J46").Value = 2312.7144
$ws.Range("K46").Value = 1611.4286
$ws.Range("L46").Value = 2312.7144
$ws.Range("M46").Value = -1423.4286
$ws.Range("N46").Value = -2688.7144
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H100").Value = 1998.25
$ws.Range("I100").Value = 1731
$ws.Range("K100").Value = 1731
$ws.Range("M100").Value = -1190
$ws.Range("H122").Value = 3563.5454
$ws.Range("I122").Value = 3333.1667
$ws.Range("J122").Value = 3840
$ws.Range("K122").Value = 9999.500100000001
$ws.Range("L122").Value = 11520
$ws.Range("M122").Value = -7549.500100000001
$ws.Range("N122").Value = -16420

# ---- Sheet: WVR (16 cell(s)) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47109
$ws.Range("J70").Value = 47109
$ws.Range("L70").Value = 47109
$ws.Range("N70").Value = -47739
$ws.Range("H73").Value = 47109
$ws.Range("J73").Value = 47109
$ws.Range("L73").Value = 47109
$ws.Range("N73").Value = -49293
$ws.Range("H113").Value = 813.5
$ws.Range("I113").Value = 719.1818
$ws.Range("K113").Value = 2157.5454
$ws.Range("M113").Value = 12.45460000000003
$ws.Range("H126").Value = 5188.5
$ws.Range("I126").Value = 5084.6665
$ws.Range("K126").Value = 15253.9995
$ws.Range("M126").Value = -12783.9995

